$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts CASA LEO FERRETERIA .. RAFAEL FERRETERIA down by one)
$ws.Rows.Item(4).Insert()

# Update clave for ROBERTO and JAIME
$ws.Range("A2").Value = 1218
$ws.Range("A3").Value = 1959

# New row 4: CHRISTIAN
$ws.Range("A4").Value = 5625
$ws.Range("B4").Value = "CHRISTIAN"
$ws.Range("D4").Value = "vendedor_estandar "

# Update rol for the rows that moved down
$ws.Range("D5").Value = "cliente_estandar"
$ws.Range("D6").Value = "cliente_estandar"
$ws.Range("D8").Value = "cliente_estandar"

# Move the selection to D4
$ws.Range("D4").Select()
